$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 486 (shifts existing rows 486:593 down to 489:596)
$ws.Rows("486:488").Insert()

# Row 486: new "Extra" record for date 2023-01-06 (serial 44932)
$ws.Cells.Item(486, 1).Value = 3
$ws.Cells.Item(486, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(486, 3).Value = "Coquimbo"
$ws.Cells.Item(486, 4).Value = 44932
$ws.Cells.Item(486, 5).Value = 5
$ws.Cells.Item(486, 6).Value = 100112028
$ws.Cells.Item(486, 7).Value = "Sandia"
$ws.Cells.Item(486, 8).Value = "Sin especificar"
$ws.Cells.Item(486, 9).Value = "Extra"
$ws.Cells.Item(486, 10).Value = 680
$ws.Cells.Item(486, 11).Value = 3800
$ws.Cells.Item(486, 12).Value = 4000
$ws.Cells.Item(486, 13).Value = 3903
$ws.Cells.Item(486, 14).Value = "`$/unidad"
$ws.Cells.Item(486, 15).Value = "Paine"
$ws.Cells.Item(486, 16).Value = 3903
$ws.Cells.Item(486, 17).Value = 1
$ws.Cells.Item(486, 18).Value = "Hortaliza"

# Row 487: new "Primera" record for date 2023-01-06 (serial 44932)
$ws.Cells.Item(487, 1).Value = 3
$ws.Cells.Item(487, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(487, 3).Value = "Coquimbo"
$ws.Cells.Item(487, 4).Value = 44932
$ws.Cells.Item(487, 5).Value = 5
$ws.Cells.Item(487, 6).Value = 100112028
$ws.Cells.Item(487, 7).Value = "Sandia"
$ws.Cells.Item(487, 8).Value = "Sin especificar"
$ws.Cells.Item(487, 9).Value = "Primera"
$ws.Cells.Item(487, 10).Value = 690
$ws.Cells.Item(487, 11).Value = 2800
$ws.Cells.Item(487, 12).Value = 3000
$ws.Cells.Item(487, 13).Value = 2904
$ws.Cells.Item(487, 14).Value = "`$/unidad"
$ws.Cells.Item(487, 15).Value = "Paine"
$ws.Cells.Item(487, 16).Value = 2904
$ws.Cells.Item(487, 17).Value = 1
$ws.Cells.Item(487, 18).Value = "Hortaliza"

# Row 488: new "Segunda" record for date 2023-01-06 (serial 44932)
$ws.Cells.Item(488, 1).Value = 3
$ws.Cells.Item(488, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(488, 3).Value = "Coquimbo"
$ws.Cells.Item(488, 4).Value = 44932
$ws.Cells.Item(488, 5).Value = 5
$ws.Cells.Item(488, 6).Value = 100112028
$ws.Cells.Item(488, 7).Value = "Sandia"
$ws.Cells.Item(488, 8).Value = "Sin especificar"
$ws.Cells.Item(488, 9).Value = "Segunda"
$ws.Cells.Item(488, 10).Value = 350
$ws.Cells.Item(488, 11).Value = 2000
$ws.Cells.Item(488, 12).Value = 2000
$ws.Cells.Item(488, 13).Value = 2000
$ws.Cells.Item(488, 14).Value = "`$/unidad"
$ws.Cells.Item(488, 15).Value = "Paine"
$ws.Cells.Item(488, 16).Value = 2000
$ws.Cells.Item(488, 17).Value = 1
$ws.Cells.Item(488, 18).Value = "Hortaliza"
